$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "Garage/Vehicle Services > Diagnostic Testing"
$ws.Range("F6").Value = "Shrewsbury"

$ws.Range("F6").Select() | Out-Null
